$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the "Hyperlink" and "Unresolved Mention" character styles to the
#    document (Word normally mints these built-in styles the first time a
#    hyperlink / @mention is inserted and they are not already present).
# ---------------------------------------------------------------------------

$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = $d.Styles("DefaultParagraphFont")
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkFont = $hyperlinkStyle.Font
$hyperlinkFont.TextColor.ObjectThemeColor = 10   # wdThemeColorHyperlink
$hyperlinkFont.Underline = 1                     # wdUnderlineSingle

$mentionStyle = $d.Styles.Add("UnresolvedMention", 2)
$mentionStyle.NameLocal = "Unresolved Mention"
$mentionStyle.BaseStyle = $d.Styles("DefaultParagraphFont")
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionFont = $mentionStyle.Font
$mentionFont.TextColor.RGB = 6053472             # 0x605E5C

# ---------------------------------------------------------------------------
# 2) Fill in the previously-empty paragraph right after "Video URL:" with the
#    Loom share link, formatted as a hyperlink, followed by a trailing space.
# ---------------------------------------------------------------------------

$findRange = $d.Content
$findRange.Find.Execute("Video URL:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$videoUrlPara = $findRange.Paragraphs(1)
$targetPara = $videoUrlPara.Next()
$targetRange = $targetPara.Range
$targetRange.Collapse(1)

$url = "https://www.loom.com/share/1a062032fe054a0e8a206487d1cb56d4?sid=22b8ceab-1aed-49b5-b78f-03c94e517169"
$link = $d.Hyperlinks.Add($targetRange, $url)

$afterLink = $d.Range($link.Range.End, $link.Range.End)
$afterLink.InsertAfter(" ")
